# Updates cryptos list price (column D) and 1h volume % (column E) values.
# Prices in column D are stored as literal text (not numbers), so each one
# is written through a temporary Text number format to stop Excel from
# auto-converting a numeric-looking string into a real number; the cell is
# then restored to the default "Normal" style so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.348.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.662.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("E8").Value = "  +1.50%  "

$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.09%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.894.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.656.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.28%  "

$ws.Range("E15").Value = "  +1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.317.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("E30").Value = "  +1.39%  "

$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "

$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.805.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0984"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
